# Reconcile henchmen name text with the printed card text, and leave the
# selection on the "data" sheet (cell A3) as the active view.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("data")

$ws.Range("A40").Value = "HYDRA Base"
$ws.Range("A8").Value  = "M.O.D.O.K.s"
$ws.Range("A16").Value = "Multiple Man"
$ws.Range("A23").Value = "Shi'ar Patrol Craft"
$ws.Range("A22").Value = "Shi'ar Death Commandos"
$ws.Range("A17").Value = "S.H.I.E.L.D. Assault Squad"
$ws.Range("A2").Value  = "Spider-Infected"

# Make "data" the active sheet/tab, with A3 selected (matches the final
# state recorded in the workbook: activeTab removed from bookViews,
# tabSelected moved from "meta" to "data", and the selection moved there).
$ws.Select()
$ws.Range("A3").Select()
